$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (new D value, new E value). $null means "leave unchanged".
$updates = @{
    2 = @("301.44", "0.22%")
    3 = @("31.42", "-0.38%")
    4 = @("5.084", "-1.24%")
    5 = @("0.07868", "-1.16%")
    6 = @("2.325", "0.37%")
    7 = @("7.813", "-1.57%")
    8 = @("3.837", "-0.54%")
    9 = @("0.9200", "0.61%")
    10 = @($null, "0.34%")
    11 = @("0.07604", "3.84%")
    12 = @("0.09067", "12.22%")
    13 = @("0.03030", "-1.33%")
    14 = @($null, "0.65%")
    15 = @("0.001519", "0.50%")
    16 = @("0.006143", "1.01%")
    17 = @("3.474", "-0.70%")
    18 = @("2.242", "0.49%")
    19 = @($null, "0.58%")
    20 = @("0.1289", "-4.73%")
    21 = @("4.122", "-12.04%")
    22 = @($null, "6.14%")
    23 = @("0.04607", "-0.47%")
    24 = @("0.001253", "-1.37%")
    25 = @("0.004466", "-0.05%")
    26 = @("0.0001247", "4.27%")
    27 = @("0.0003385", "-1.92%")
    39 = @("0.01751", "-4.70%")
    40 = @("0.04637", "2.07%")
    41 = @("0.007017", "-4.02%")
    42 = @("0.1359", "1.04%")
    43 = @("0.002185", "0.40%")
    44 = @("0.009757", "-8.30%")
    45 = @("0.00006258", "-2.43%")
    46 = @("0.00000000749", "-0.64%")
    47 = @($null, "19.49%")
    48 = @("1.153", "40.51%")
    49 = @("0.00002098", "-0.64%")
    50 = @("0.0001998", "-0.64%")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]

    # Build a single contiguous range address covering the cell(s) that change
    # on this row, e.g. "D2:E2", "E10:E10" etc. Using one contiguous Range object
    # (rather than a comma-joined multi-area address) ensures the text-number-format
    # trick below is applied consistently to every touched cell.
    if (($dVal -ne $null) -and ($eVal -ne $null)) {
        $rangeAddr = "D" + $row + ":E" + $row
    } elseif ($dVal -ne $null) {
        $rangeAddr = "D" + $row + ":D" + $row
    } else {
        $rangeAddr = "E" + $row + ":E" + $row
    }
    $rng = $ws.Range($rangeAddr)

    # Force the cells to be stored as text (matching the original inline-string
    # cells) instead of letting Excel auto-convert numeric-looking strings into
    # numbers or percentages.
    $rng.NumberFormat = "@"

    if ($dVal -ne $null) { $ws.Range("D$row").Value = $dVal }
    if ($eVal -ne $null) { $ws.Range("E$row").Value = $eVal }

    # Restore default (General) formatting/style so the cells look the same as
    # they did before (no special number format or style index).
    $rng.NumberFormatLocal = "General"
    $rng.Style = "Normal"
}
